# Tekken 8 Steve Fox frame data: convert the Block (E) and Hit (F) columns
# from free-form text (e.g. "+1", "-3", "KND", "Launch", "+17/KND") into
# plain numeric values.
#
#   Block (E): the textual value is a signed number (optionally with a
#     trailing letter like "c"/"g", or malformed like "-+2"). The stored
#     sign convention is inverted compared to the numeric column, so the
#     numeric result is the NEGATION of the parsed signed number
#     (e.g. "+1" -> -1, "-3" -> 3). "KND" or unparsable text -> blank.
#
#   Hit (F): the textual value's leading signed number is kept as-is
#     (sign preserved, e.g. "+8" -> 8, "-2" -> -2). "KND" / "Launch" (and
#     combined forms like "+17/KND") collapse to 70. Unparsable text that
#     has no leading number (e.g. "+KND") -> blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Convert-Block {
    param([string]$text)

    if ([string]::IsNullOrEmpty($text)) { return $null }
    $t = $text.Trim()
    if ($t.Length -eq 0) { return $null }
    if ($t.ToUpper() -eq "KND") { return $null }
    if ($t -eq "0") { return 0 }
    if ($t -match '^([+-])(\d+)[A-Za-z]*$') {
        $sign = $matches[1]
        $num = [int]$matches[2]
        if ($sign -eq '-') { return $num } else { return (0 - $num) }
    }
    return $null
}

function Convert-Hit {
    param([string]$text)

    if ([string]::IsNullOrEmpty($text)) { return $null }
    $t = $text.Trim()
    if ($t.Length -eq 0) { return $null }
    $u = $t.ToUpper()
    if ($u -eq "KND" -or $u -eq "LAUNCH") { return 70 }
    if ($t -match '^([+-]?)(\d+)') {
        $sign = $matches[1]
        $num = [int]$matches[2]
        if ($sign -eq '-') { return (0 - $num) } else { return $num }
    }
    return $null
}

$lastRow = $ws.Range("A1").End(4).Row   # xlDown
if ($lastRow -lt 2) { $lastRow = 154 }

for ($r = 2; $r -le $lastRow; $r++) {
    $eCell = $ws.Range("E$r")
    $eOld = $eCell.Value2
    if (-not [string]::IsNullOrEmpty($eOld)) {
        $eNew = Convert-Block([string]$eOld)
        if ($null -eq $eNew) {
            $eCell.Value = ""
        } else {
            $eCell.Value = $eNew
        }
    }

    $fCell = $ws.Range("F$r")
    $fOld = $fCell.Value2
    if (-not [string]::IsNullOrEmpty($fOld)) {
        $fNew = Convert-Hit([string]$fOld)
        if ($null -eq $fNew) {
            $fCell.Value = ""
        } else {
            $fCell.Value = $fNew
        }
    }
}
